## Insert a new weekly price record for "Jengibre" (Vega Modelo de Temuco)
## at row 35 - this pushes the existing rows 35..93 down to 36..94 and
## appends the former last row's data as the new row 94 (handled
## automatically by the row insert / shift).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert one blank row above the current row 35, shifting row 35..93 (and
# all their formatting) down to 36..94.
$ws.Rows(35).Insert()

# Populate the newly inserted row 35 with the new observation.
$ws.Cells.Item(35, 1).Value  = 10
$ws.Cells.Item(35, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(35, 3).Value  = "La Araucanía"
$ws.Cells.Item(35, 4).Value  = 44413
$ws.Cells.Item(35, 5).Value  = 9
$ws.Cells.Item(35, 6).Value  = 100114007
$ws.Cells.Item(35, 7).Value  = "Jengibre"
$ws.Cells.Item(35, 8).Value  = "Sin especificar"
$ws.Cells.Item(35, 9).Value  = "Primera"
$ws.Cells.Item(35, 10).Value = 100
$ws.Cells.Item(35, 11).Value = 25000
$ws.Cells.Item(35, 12).Value = 25000
$ws.Cells.Item(35, 13).Value = 25000
$ws.Cells.Item(35, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(35, 15).Value = "Perú"
$ws.Cells.Item(35, 16).Value = 1923
$ws.Cells.Item(35, 17).Value = 13
$ws.Cells.Item(35, 18).Value = "Hortaliza"
